# cryptos.xlsx refresh (GitHub Actions bot) -- rewrite the Price (D) and
# Volume(1h) (E) columns for every coin row with the latest scraped
# figures, and apply the FraxShare/VeChain rank swap at rows 42-43
# (new coin, link, price and volume for each of those two rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings in the source data (several
# are already non-numeric, e.g. "26.664.09"). Cells whose new price
# text parses as a plain number need a leading quote-prefix apostrophe
# so Excel keeps storing them as TEXT instead of silently converting
# them to a floating point number.

$ws.Range("D2").Value = '26.664.09'
$ws.Range("E2").Value = '  +4.33%  '
$ws.Range("D3").Value = '1.747.62'
$ws.Range("E3").Value = '  +4.59%  '
$ws.Range("D4").Value = '''0.9992'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''247.25'
$ws.Range("E5").Value = '  +3.42%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '''0.4807'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").Value = '''0.2713'
$ws.Range("E8").Value = '  +3.03%  '
$ws.Range("D9").Value = '''0.06247'
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").Value = '1.747.03'
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").Value = '''0.07113'
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '''15.84'
$ws.Range("E12").Value = '  +6.43%  '
$ws.Range("D13").Value = '''0.6195'
$ws.Range("E13").Value = '  +4.68%  '
$ws.Range("D14").Value = '''4.512'
$ws.Range("E14").Value = '  +2.79%  '
$ws.Range("D15").Value = '''77.31'
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").Value = '''0.9989'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '26.673.36'
$ws.Range("E17").Value = '  +4.41%  '
$ws.Range("D18").Value = '''0.9994'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '''0.000006895'
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").Value = '''11.73'
$ws.Range("E20").Value = '  +2.07%  '
$ws.Range("D21").Value = '1.970.14'
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("D22").Value = '''4.650'
$ws.Range("E22").Value = '  +4.60%  '
$ws.Range("D23").Value = '''8.879'
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("D24").Value = '''5.361'
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("D25").Value = '''136.44'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '''15.49'
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("D27").Value = '''1.830'
$ws.Range("E27").Value = '  +5.92%  '
$ws.Range("D28").Value = '''1.417'
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("D29").Value = '''107.81'
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("D30").Value = '''4.033'
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").Value = '''3.768'
$ws.Range("E31").Value = '  +3.04%  '
$ws.Range("D32").Value = '''0.07897'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").Value = '''0.04580'
$ws.Range("E33").Value = '  +8.48%  '
$ws.Range("D34").Value = '''2.612'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = '''0.6395'
$ws.Range("E35").Value = '  +4.68%  '
$ws.Range("D36").Value = '''0.9985'
$ws.Range("E36").Value = '  +4.92%  '
$ws.Range("D37").Value = '''0.9508'
$ws.Range("E37").Value = '  +10.97%  '
$ws.Range("D38").Value = '''114.19'
$ws.Range("E38").Value = '  +18.67%  '
$ws.Range("D39").Value = '''2.482'
$ws.Range("E39").Value = '  -4.56%  '
$ws.Range("D40").Value = '''1.980'
$ws.Range("E40").Value = '  +5.61%  '
$ws.Range("D41").Value = '''1.003'
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.01521'
$ws.Range("E42").Value = '  +3.11%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.741'
$ws.Range("E43").Value = '  +18.09%  '
$ws.Range("D44").Value = '''0.3923'
$ws.Range("E44").Value = '  +3.77%  '
$ws.Range("D45").Value = '''6.729'
$ws.Range("E45").Value = '  +8.00%  '
$ws.Range("D46").Value = '''0.1205'
$ws.Range("E46").Value = '  +7.82%  '
$ws.Range("D47").Value = '''0.05333'
$ws.Range("E47").Value = '  +1.36%  '
$ws.Range("D48").Value = '''8.028'
$ws.Range("E48").Value = '  +8.60%  '
$ws.Range("D49").Value = '''30.89'
$ws.Range("E49").Value = '  +3.43%  '
$ws.Range("D50").Value = '''0.3459'
$ws.Range("E50").Value = '  +3.38%  '
$ws.Range("D51").Value = '''51.76'
$ws.Range("E51").Value = '  +3.28%  '
